$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a literal text value (avoids Excel auto-converting
    # strings like "58%" into a numeric percentage).
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# Row 2
$ws.Range("C2").Value = "25°"
Set-TextValue "D2" "58%"

# Row 3
$ws.Range("C3").Value = "25°"
Set-TextValue "D3" "60%"
Set-TextValue "E3" "84%"

# Row 4
$ws.Range("B4").Value = "35°"
$ws.Range("C4").Value = "25°"
Set-TextValue "D4" "62%"
Set-TextValue "E4" "84%"

# Row 5
$ws.Range("B5").Value = "33°"
$ws.Range("C5").Value = "25°"
Set-TextValue "D5" "67%"
Set-TextValue "E5" "87%"

# Row 6
Set-TextValue "D6" "72%"
Set-TextValue "E6" "89%"

# Row 7
Set-TextValue "E7" "88%"

# Row 8
$ws.Range("C8").Value = "24°"
Set-TextValue "D8" "76%"
Set-TextValue "E8" "93%"

# Row 9
$ws.Range("B9").Value = "33°"
$ws.Range("C9").Value = "24°"
Set-TextValue "D9" "72%"
Set-TextValue "E9" "93%"

# Row 10
$ws.Range("B10").Value = "32°"
$ws.Range("C10").Value = "24°"
Set-TextValue "D10" "73%"
Set-TextValue "E10" "94%"

# Row 11
$ws.Range("B11").Value = "32°"
Set-TextValue "D11" "72%"
Set-TextValue "E11" "90%"
